# Update forest data - 2025-12-03 12:20
#
# Workflow: the listings most recently scraped live on the "New" sheet.
# On every refresh those listings are archived onto the bottom of
# "Previously added", and "New" is repopulated with the newest batch.

$wb     = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ----------------------------------------------------------------------
# Step 1 - archive the rows currently on "New" (rows 2-4) onto the end
# of "Previously added" (currently rows 2-313, landing on 314-316).
# ----------------------------------------------------------------------
$prevLastRow = $wsPrev.UsedRange.Rows.Count   # 313
$newLastRow  = $wsNew.UsedRange.Rows.Count    # 4

$archiveFirstRow = $prevLastRow + 1

$srcRange = $wsNew.Range("A2:F" + $newLastRow)
$srcRange.Copy($wsPrev.Cells.Item($archiveFirstRow, 1))

# Links that were on "New" and are now archived.
$archivedLinks = @(
    "https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/aluksne/gxkjp.html",
    "https://www.ss.com/msg/lv/real-estate/wood/ventspils-and-reg/zleku-pag/booep.html",
    "https://www.ss.com/msg/lv/real-estate/wood/ventspils-and-reg/zleku-pag/bbkhx.html"
)

for ($i = 0; $i -lt $archivedLinks.Count; $i++) {
    $destRow = $archiveFirstRow + $i
    $destCell = $wsPrev.Cells.Item($destRow, 1)
    $wsPrev.Hyperlinks.Add($destCell, $archivedLinks[$i]) | Out-Null
    # Adding a hyperlink re-styles the cell with the generic "Hyperlink"
    # look; restore the sheet's own link style used by every other row.
    $wsPrev.Range("A2").Copy()
    $destCell.PasteSpecial(-4122) | Out-Null
}

# ----------------------------------------------------------------------
# Step 2 - replace "New" with the freshly scraped listings (5 rows).
# Extend the existing per-column formatting down to rows 5 & 6 first.
# ----------------------------------------------------------------------
$formatSrc = $wsNew.Range("A4:F4")
$formatSrc.Copy($wsNew.Range("A5:F5"))
$formatSrc.Copy($wsNew.Range("A6:F6"))

# Drop the 3 old hyperlinks; they'll be rebuilt for the new rows below.
$wsNew.Range("A2:F" + $newLastRow).Hyperlinks.Delete()

$newRows = @(
    @{ link = "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/cghpfg.html"; price = "8 000 €";   district = "Jelgava un raj."; area = "2 ha.";  cadastre = "54860020101";  date = 45994.59513888889 },
    @{ link = "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/livberzes-pag/ohbhg.html";  price = "4 500 €";   district = "Jelgava un raj."; area = "1 ha.";  cadastre = "54620090154";  date = 45994.58333333333 },
    @{ link = "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/egcjx.html";  price = "4 000 €";   district = "Jelgava un raj."; area = "1 ha.";  cadastre = "548600701105"; date = 45994.57708333334 },
    @{ link = "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/ecmgi.html";  price = "35 000 €";  district = "Jelgava un raj."; area = "10 ha."; cadastre = "54860060066";  date = 45994.57152777778 },
    @{ link = "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/varmes-pag/dpmoh.html";     price = "149 000 €"; district = "Kuldīga un raj."; area = "38 ha."; cadastre = "";             date = 45993.70833333333 }
)

# Untouched style templates to restore from after styling gets stomped
# on by Hyperlinks.Add / NumberFormat below. Row 2's A/E cells are never
# themselves used as a paste target, so they stay a clean reference.
$aStyleTemplate = $wsNew.Range("A2")
$eStyleTemplate = $wsNew.Range("E2")

$row = 2
foreach ($item in $newRows) {
    $aCell = $wsNew.Cells.Item($row, 1)
    $eCell = $wsNew.Cells.Item($row, 5)

    $aCell.Value2 = $item.link
    $wsNew.Cells.Item($row, 2).Value2 = $item.price
    $wsNew.Cells.Item($row, 3).Value2 = $item.district
    $wsNew.Cells.Item($row, 4).Value2 = $item.area

    # Cadastre numbers are all-digit strings; force text so they keep
    # being stored (and shared) as plain text, not coerced to numbers.
    $eCell.NumberFormat = "@"
    $eCell.Value2 = $item.cadastre

    $wsNew.Cells.Item($row, 6).Value2 = $item.date

    $wsNew.Hyperlinks.Add($aCell, $item.link) | Out-Null

    # Restore this row's original per-column styles (hyperlink-add and
    # the text numberformat tweak both stomp on cell styling).
    $aStyleTemplate.Copy()
    $aCell.PasteSpecial(-4122) | Out-Null
    $eStyleTemplate.Copy()
    $eCell.PasteSpecial(-4122) | Out-Null

    $row++
}

Write-Output "forest data updated"
